$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently shows "R40" (shared text). The target value is the literal
# text "1" -- a numeric-looking string that must stay text (not become the
# number 1). Writing it directly via .Value auto-converts numeric-looking
# strings to real numbers, so stage it as text in a scratch cell (using the
# leading-apostrophe "force text" entry) and paste just the *value* into
# B11 so the destination's existing formatting/style is left untouched.
$scratch = $ws.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
